# Update 20 - extent reports 3.0 with screen shot
#
# On the "Test" worksheet, the Runmode values for the two LoginSuite test
# cases (rows 14 and 15) are swapped: TC01 now runs ("Y") and TC02 is
# turned off ("N").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

$tc01Runmode = $ws.Range("B14").Value2
$tc02Runmode = $ws.Range("B15").Value2

$ws.Range("B14").Value2 = $tc02Runmode
$ws.Range("B15").Value2 = $tc01Runmode

$ws.Activate()
$ws.Range("B15").Select()
